$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove extra rows 22-26 (5 rows) so only rows 1-21 remain
$ws.Range("A22:C26").EntireRow.Delete()

# Update data rows 2-21 with new article data
$ws.Range("A2").Value = "Hàng nghìn biệt thự bỏ hoang nơi đất đấu giá 130 triệu một m2"
$ws.Range("B2").Value = "https://vnexpress.net/hang-nghin-biet-thu-bo-hoang-noi-dat-dau-gia-130-trieu-mot-m2-4784776.html"
$ws.Range("C2").Value = "Hà NộiThửa đất nền ""chưa có gì"" ở Hoài Đức được tranh mua xuyên đêm, giá hơn 130 triệu đồng trong khi nhà liền kề, biệt thự xây sẵn gần đó không người ở."

$ws.Range("A3").Value = "Trước khi tôi mất tiền…"
$ws.Range("B3").Value = "https://vnexpress.net/truoc-khi-toi-mat-tien-4785171.html"
$ws.Range("C3").Value = "Đang lướt Facebook, mắt tôi vấp phải đường link bài viết có tiêu đề gây sốc: 'Ngân hàng Nhà nước kiện ông A'."

$ws.Range("A4").Value = "Đề xuất xây dựng một đoạn cao tốc Hà Nội - Viêng Chăn"
$ws.Range("B4").Value = "https://vnexpress.net/de-xuat-xay-dung-mot-doan-cao-toc-ha-noi-vieng-chan-4785334.html"
$ws.Range("C4").Value = "Dự án cao tốc Hà Nội - Viêng Chăn, đoạn Vinh - Thanh Thủy (Nghệ An) dài 65 km được đề xuất xây dựng với tổng vốn 18.500 tỷ đồng."

$ws.Range("A5").Value = "Thông điệp từ chuyến thăm Ukraine của Thủ tướng Ấn Độ"
$ws.Range("B5").Value = "https://vnexpress.net/thong-diep-tu-chuyen-tham-ukraine-cua-thu-tuong-an-do-4785116.html"
$ws.Range("C5").Value = "Lần đầu thăm Ukraine, Thủ tướng Modi dường như muốn trấn an phương Tây rằng Ấn Độ không hoàn toàn ngả về Nga như họ vẫn nghĩ."

$ws.Range("A6").Value = "Cổ vật Hoàng thành Thăng Long trưng bày ở TP HCM"
$ws.Range("B6").Value = "https://vnexpress.net/co-vat-hoang-thanh-thang-long-trung-bay-o-tp-hcm-4784829.html"
$ws.Range("C6").Value = "150 hiện vật, tài liệu, hình ảnh tại Khu di sản Hoàng thành Thăng Long trưng bày ở Bảo tàng TP HCM, quận 1."

$ws.Range("A7").Value = "Arsenal đòi được món nợ từ Aston Villa"
$ws.Range("B7").Value = "https://vnexpress.net/arsenal-doi-duoc-mon-no-tu-aston-villa-4785356.html"
$ws.Range("C7").Value = "AnhTận dụng cơ hội tốt hơn kèm một chút may mắn, thầy trò Mikel Arteta đánh bại đối thủ khó chơi 2-0 ở vòng 2 Ngoại hạng Anh."

$ws.Range("A8").Value = "Bộ Giáo dục: Thầy cô được đàng hoàng dạy thêm"
$ws.Range("B8").Value = "https://vnexpress.net/bo-giao-duc-thay-co-duoc-dang-hoang-day-them-4785315.html"
$ws.Range("C8").Value = "Thầy cô được đàng hoàng dạy học sinh của mình ngoài nhà trường nhưng tuyệt đối không được ép buộc, theo Vụ trưởng Trung học."

$ws.Range("A9").Value = "IS nhận trách nhiệm vụ đâm dao tại lễ hội ở Đức"
$ws.Range("B9").Value = "https://vnexpress.net/is-nhan-trach-nhiem-vu-dam-dao-tai-le-hoi-o-duc-4785362.html"
$ws.Range("C9").Value = "IS tuyên bố kẻ thực hiện vụ đâm dao khiến ba người chết tại lễ hội ở thành phố Solingen, phía tây Đức, là thành viên của nhóm này."

$ws.Range("A10").Value = "iPhone 16 có thể bán tại Việt Nam cuối tháng 9"
$ws.Range("B10").Value = "https://vnexpress.net/iphone-16-co-the-ban-tai-viet-nam-cuoi-thang-9-4784978.html"
$ws.Range("C10").Value = "Việt Nam vẫn nằm trong những thị trường ưu tiên thứ hai của Apple và nhiều khả năng sẽ bán iPhone 16 ngay trong tháng 9."

$ws.Range("A11").Value = "Hà Nội phân luồng giao thông dịp nghỉ lễ 2/9"
$ws.Range("B11").Value = "https://vnexpress.net/ha-noi-phan-luong-giao-thong-dip-nghi-le-2-9-4785331.html"
$ws.Range("C11").Value = "Nhằm hạn chế ùn tắc tại các tuyến đường ra, vào thành phố dịp lễ 2/9, Sở Giao thông Vận tải Hà Nội vừa thông báo hướng dẫn phân luồng."

$ws.Range("A12").Value = "Tiến sĩ hóa học 'bước ra' từ gian bếp nghèo của mẹ"
$ws.Range("B12").Value = "https://vnexpress.net/tien-si-hoa-hoc-buoc-ra-tu-gian-bep-ngheo-cua-me-4784851.html"
$ws.Range("C12").Value = "Hà NộiNăm 2012, khi nhận bằng tiến sĩ loại xuất sắc, các phóng viên Tây Ban Nha hỏi lý do đến đây học, cô gái Việt Nam Vũ Thị Tần liền bắt đầu câu chuyện từ gian bếp của mẹ."

$ws.Range("A13").Value = "Haaland lập hat-trick giúp Man City thắng ngược"
$ws.Range("B13").Value = "https://vnexpress.net/haaland-lap-hat-trick-giup-man-city-thang-nguoc-4785344.html"
$ws.Range("C13").Value = "AnhTiền đạo Erling Haaland lập hat-trick, giúp chủ nhà Man City đè bẹp tân binh Ipswich Town 4-1 ở vòng hai Ngoại hạng Anh."

$ws.Range("A14").Value = "Chuyến ăn mừng trên du thuyền hóa thảm kịch của tỷ phú Anh"
$ws.Range("B14").Value = "https://vnexpress.net/chuyen-an-mung-tren-du-thuyen-hoa-tham-kich-cua-ty-phu-anh-4785120.html"
$ws.Range("C14").Value = "Tỷ phú Mike Lynch muốn đi du thuyền khắp Địa Trung Hải để ăn mừng phán quyết vô tội tại Mỹ, nhưng hành trình nhanh chóng biến thành thảm kịch."

$ws.Range("A15").Value = "Tường San đoạt á hậu Chuyển giới Quốc tế"
$ws.Range("B15").Value = "https://vnexpress.net/tuong-san-doat-a-hau-chuyen-gioi-quoc-te-4785306.html"
$ws.Range("C15").Value = "Thái LanTường San, 19 tuổi, đoạt danh hiệu á hậu 2 ở chung kết Hoa hậu Chuyển giới Quốc tế lần 18, tối 24/8."

$ws.Range("A16").Value = "Thanh niên chết não hiến tạng ghép cho 6 người"
$ws.Range("B16").Value = "https://vnexpress.net/thanh-nien-chet-nao-hien-tang-ghep-cho-6-nguoi-4785312.html"
$ws.Range("C16").Value = "Hà NộiNam thanh niên 32 tuổi chết não do tai nạn giao thông, hiến hai quả thận, giác mạc, tim, gan, là trường hợp lấy - ghép mô tạng đầu tiên do Bệnh viện Đa khoa Xanh Pôn thực hiện."

$ws.Range("A17").Value = "Son Heung-min ghi cú đúp khi Tottenham thắng đậm"
$ws.Range("B17").Value = "https://vnexpress.net/son-heung-min-ghi-cu-dup-khi-tottenham-thang-dam-4785351.html"
$ws.Range("C17").Value = "AnhTiền đạo Hàn Quốc Son Heung-min lập cú đúp, giúp Tottenham thắng đội khách Everton 4-0 ở vòng hai Ngoại hạng Anh."

$ws.Range("A18").Value = "Cá nhân, chủ hộ kinh doanh nợ thuế có thể bị cấm xuất cảnh"
$ws.Range("B18").Value = "https://vnexpress.net/ca-nhan-chu-ho-kinh-doanh-no-thue-co-the-bi-cam-xuat-canh-4785309.html"
$ws.Range("C18").Value = "Bộ Tài chính muốn thêm cá nhân, chủ hộ kinh doanh vào đối tượng bị tạm hoãn xuất cảnh do chưa hoàn thành nghĩa vụ thuế."

$ws.Range("A19").Value = "Iran muốn kiểm soát tình trạng thù địch với Mỹ"
$ws.Range("B19").Value = "https://vnexpress.net/iran-muon-kiem-soat-tinh-trang-thu-dich-voi-my-4785338.html"
$ws.Range("C19").Value = "Ngoại trưởng Iran Abbas Araghchi tuyên bố nước này muốn kiểm soát tình trạng thù địch với Mỹ để giảm bớt sức ép và đối phó lệnh trừng phạt."

$ws.Range("A20").Value = "Chuyên gia: Vị thế của TP HCM đang bị 'xói mòn'"
$ws.Range("B20").Value = "https://vnexpress.net/chuyen-gia-vi-the-cua-tp-hcm-dang-bi-xoi-mon-4785273.html"
$ws.Range("C20").Value = "Vị thế của TP HCM đang bị ""xói mòn"" và đứng trước nhiều thách thức khi so với các địa phương khác và một số đô thị lớn ở Đông Nam Á, theo TS Vũ Thành Tự Anh."

$ws.Range("A21").Value = "Man Utd thua phút 95 ở Ngoại hạng Anh"
$ws.Range("B21").Value = "https://vnexpress.net/man-utd-thua-phut-95-o-ngoai-hang-anh-4785327.html"
$ws.Range("C21").Value = "AnhMan Utd thất bại 1-2 trước chủ nhà Brighton với bàn thua ở phút bù hiệp hai, trận sớm nhất vòng hai Ngoại hạng Anh."
